# Apply cryptos list price/volume updates (GitHub Actions scrape refresh)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "60.242.67"
$ws.Range("E2").Value = "  -4.33%  "
$ws.Range("D3").Value = "2.984.42"
$ws.Range("E3").Value = "  -5.91%  "
$ws.Range("E4").Value = "  +0.07%  "
$ws.Range("D5").Value = "'578.98"
$ws.Range("E5").Value = "  -2.17%  "
$ws.Range("D6").Value = "'125.40"
$ws.Range("E6").Value = "  -6.84%  "
$ws.Range("E7").Value = "  +0.09%  "
$ws.Range("D8").Value = "2.981.80"
$ws.Range("E8").Value = "  -5.91%  "
$ws.Range("D9").Value = "'0.500"
$ws.Range("E9").Value = "  -2.64%  "
$ws.Range("E10").Value = "  -5.74%  "
$ws.Range("E11").Value = "  -2.20%  "
$ws.Range("D12").Value = "'0.440"
$ws.Range("E12").Value = "  -2.68%  "
$ws.Range("E13").Value = "  -5.81%  "
$ws.Range("D14").Value = "'32.51"
$ws.Range("E14").Value = "  -5.14%  "
$ws.Range("E15").Value = "  +0.32%  "
$ws.Range("D16").Value = "3.473.48"
$ws.Range("E16").Value = "  -5.94%  "
$ws.Range("D17").Value = "2.984.37"
$ws.Range("E17").Value = "  -5.85%  "
$ws.Range("D18").Value = "60.135.04"
$ws.Range("E18").Value = "  -4.46%  "
$ws.Range("D19").Value = "'6.21"
$ws.Range("E19").Value = "  -4.75%  "
$ws.Range("D20").Value = "'431.83"
$ws.Range("E20").Value = "  -6.04%  "
$ws.Range("D21").Value = "'13.10"
$ws.Range("E21").Value = "  -6.26%  "
$ws.Range("D22").Value = "'0.661"
$ws.Range("E22").Value = "  -4.87%  "
$ws.Range("E23").Value = "  -7.32%  "
$ws.Range("D24").Value = "'12.66"
$ws.Range("E24").Value = "  -4.61%  "
$ws.Range("D25").Value = "'79.04"
$ws.Range("E25").Value = "  -3.93%  "
$ws.Range("E26").Value = "  +0.12%  "
$ws.Range("E27").Value = "  -0.03%  "
$ws.Range("D29").Value = "'7.27"
$ws.Range("E29").Value = "  -4.58%  "
$ws.Range("E30").Value = "  -6.73%  "
$ws.Range("D31").Value = "'6.13"
$ws.Range("E31").Value = "  -8.79%  "
$ws.Range("E32").Value = "  -6.68%  "
$ws.Range("D33").Value = "'0.0934"
$ws.Range("E33").Value = "  -7.69%  "
$ws.Range("D34").Value = "'2.16"
$ws.Range("E34").Value = "  -8.50%  "
$ws.Range("D35").Value = "'0.952"
$ws.Range("E35").Value = "  -6.99%  "
$ws.Range("D36").Value = "'5.59"
$ws.Range("E36").Value = "  -3.28%  "
$ws.Range("D37").Value = "'49.63"
$ws.Range("E37").Value = "  -2.96%  "
$ws.Range("D38").Value = "0.0₃0658"
$ws.Range("E38").Value = "  -6.73%  "
$ws.Range("D39").Value = "'8.00"
$ws.Range("E39").Value = "  -0.90%  "
$ws.Range("E40").Value = "  -6.88%  "
$ws.Range("D41").Value = "'385.10"
$ws.Range("E41").Value = "  -4.34%  "
$ws.Range("E42").Value = "  -2.23%  "
$ws.Range("D43").Value = "'2.45"
$ws.Range("E43").Value = "  -7.04%  "
$ws.Range("D44").Value = "2.630.28"
$ws.Range("D47").Value = "'119.49"
$ws.Range("E47").Value = "  -3.53%  "
$ws.Range("E48").Value = "  -5.43%  "
$ws.Range("E49").Value = "  -3.43%  "
$ws.Range("D50").Value = "'23.48"
$ws.Range("E50").Value = "  -6.28%  "
$ws.Range("D51").Value = "'31.10"
$ws.Range("E51").Value = "  -10.10%  "
